# Commit: "remove double-click-to-delete in favour of delete button"
#
# This script updates the tpivot phase-3 tracking sheet:
#  - Item #16 (double-click delete UX) is marked Complete with a finish
#    date, its Approach note gains a follow-up sentence, and its
#    Est Difficulty drops from MED to LOW.
#  - Item #11's Approach note (column/reducer aliasing) gains a follow-up
#    sentence and its Est Difficulty is bumped from MED to HIGH.
#  - Item #14's Approach note gains a "See #11" cross reference.
#  - The "changes currently modify DOM..." note is promoted to be the
#    Feature description of its own row, with a "See #11" cross
#    reference appended to the original row's Approach text.
#  - The used range / AutoFilter / _FilterDatabase grow from row 23 to
#    row 26 to cover the full table.
#  - The view scrolls back to the top and selects E7 instead of E23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 7 (Feature #11): bump difficulty, extend approach note ---
$ws.Range("D7").Value = "HIGH"
$e7 = $ws.Range("E7").Value()
$ws.Range("E7").Value = $e7 + " Perhaps pivot data should be represented as an array of cells that know their own coordinates, as in reference pivot library."
$ws.Rows.Item(7).RowHeight = 63.75

# --- Row 17 (Feature #14): add cross reference to approach note ---
$e17 = $ws.Range("E17").Value()
$ws.Range("E17").Value = $e17 + " See #11"

# --- Row 21 (Feature #16, double-click delete): completed now ---
$ws.Range("D21").Value = "LOW"
$e21 = $ws.Range("E21").Value()
$ws.Range("E21").Value = $e21 + " Use Jquery hover()"
$ws.Range("F21").Value = "Complete"
$ws.Range("G21").Value = 42978
$ws.Range("H21").Value = 42978
$ws.Rows.Item(21).RowHeight = 25.5

# --- Row 24: feature text replaced, approach note gains cross reference ---
$ws.Range("B24").Value = "Changes currently modify DOM, not app state. Conversely, charts are generated from state (specifically most recently received pivot results). Need to pin shape modifications to state."
$e24 = $ws.Range("E24").Value()
$ws.Range("E24").Value = $e24 + " See #11"

# --- Grow the table's filter / used range from row 23 to row 26 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:H26").AutoFilter()
$fd = $wb.Names.Item("Sheet1!_FilterDatabase")
$fd.RefersTo = "=Sheet1!`$A`$1:`$H`$26"

# --- Reset the view: scroll to top, select E7 ---
$wnd = $excel.ActiveWindow
$wnd.ScrollRow = 1
$wnd.ScrollColumn = 1
$ws.Range("E7").Select()
